$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values (Id, coordinates, taxon info, activity/comment fields)
$ws.Range("A2").Value2 = 104449298
$ws.Range("B2").Value2 = 56395
$ws.Range("E2").Value2 = 100109
$ws.Range("F2").Value2 = "Tretåig hackspett"
$ws.Range("G2").Value2 = "Picoides tridactylus"
$ws.Range("H2").Value2 = "(Linnaeus, 1758)"
$ws.Range("K2").Value2 = ""
$ws.Range("L2").Value2 = ""
$ws.Range("M2").Value2 = "äldre spår"
$ws.Range("N2").Value2 = ""
$ws.Range("Q2").Value2 = 418152.1433075544
$ws.Range("R2").Value2 = 7018755.833866266
$ws.Range("AC2").Value2 = "ringhack"
$ws.Range("A3").Value2 = 104449306
$ws.Range("Q3").Value2 = 418116.3208070688
$ws.Range("R3").Value2 = 7018906.624424814
$ws.Range("A4").Value2 = 104449305
$ws.Range("M4").Value2 = "färska spår"
$ws.Range("Q4").Value2 = 418113.1107625436
$ws.Range("R4").Value2 = 7018904.455793464
$ws.Range("A5").Value2 = 104449307
$ws.Range("M5").Value2 = "äldre spår"
$ws.Range("Q5").Value2 = 418106.0762497109
$ws.Range("R5").Value2 = 7018911.38607322
$ws.Range("A6").Value2 = 104449297
$ws.Range("Q6").Value2 = 418163.1633477406
$ws.Range("R6").Value2 = 7018746.101364438
$ws.Range("A7").Value2 = 104449308
$ws.Range("M7").Value2 = "färska spår"
$ws.Range("Q7").Value2 = 418207.1051796933
$ws.Range("R7").Value2 = 7019144.644948276
$ws.Range("AC7").Value2 = "Påbörjat bo?"
$ws.Range("A8").Value2 = 104449333
$ws.Range("Q8").Value2 = 418224.807980529
$ws.Range("R8").Value2 = 7018298.906277624
$ws.Range("A9").Value2 = 104449393
$ws.Range("B9").Value2 = 78570
$ws.Range("E9").Value2 = 2081
$ws.Range("F9").Value2 = "Skrovellav"
$ws.Range("G9").Value2 = "Lobaria scrobiculata"
$ws.Range("H9").Value2 = "(Scop.) DC."
$ws.Range("Q9").Value2 = 418188.408122587
$ws.Range("R9").Value2 = 7018072.943679515

# Clear cells that should no longer have a value
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("AC9").ClearContents()
